$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update F column timestamps on the data sheet (rows 2-86)
$newTimes = @(
    "2021-10-05 14:20:01.014757",
    "2021-10-05 14:20:01.014765",
    "2021-10-05 14:20:01.014769",
    "2021-10-05 14:20:01.014771",
    "2021-10-05 14:20:01.014774",
    "2021-10-05 14:20:01.014777",
    "2021-10-05 14:20:01.014779",
    "2021-10-05 14:20:01.014782",
    "2021-10-05 14:20:01.014785",
    "2021-10-05 14:20:01.014787",
    "2021-10-05 14:20:01.014790",
    "2021-10-05 14:20:01.014792",
    "2021-10-05 14:20:01.014795",
    "2021-10-05 14:20:01.014798",
    "2021-10-05 14:20:01.014800",
    "2021-10-05 14:20:01.014803",
    "2021-10-05 14:20:01.014806",
    "2021-10-05 14:20:01.014809",
    "2021-10-05 14:20:01.014811",
    "2021-10-05 14:20:01.014814",
    "2021-10-05 14:20:01.014817",
    "2021-10-05 14:20:01.014819",
    "2021-10-05 14:20:01.014822",
    "2021-10-05 14:20:01.014824",
    "2021-10-05 14:20:01.014827",
    "2021-10-05 14:20:01.014829",
    "2021-10-05 14:20:01.014832",
    "2021-10-05 14:20:01.014834",
    "2021-10-05 14:20:01.014837",
    "2021-10-05 14:20:01.014840",
    "2021-10-05 14:20:01.014842",
    "2021-10-05 14:20:01.014845",
    "2021-10-05 14:20:01.014847",
    "2021-10-05 14:20:01.014850",
    "2021-10-05 14:20:01.014853",
    "2021-10-05 14:20:01.014855",
    "2021-10-05 14:20:01.014858",
    "2021-10-05 14:20:01.014860",
    "2021-10-05 14:20:01.014863",
    "2021-10-05 14:20:01.014865",
    "2021-10-05 14:20:01.014868",
    "2021-10-05 14:20:01.014871",
    "2021-10-05 14:20:01.014873",
    "2021-10-05 14:20:01.014876",
    "2021-10-05 14:20:01.014878",
    "2021-10-05 14:20:01.014881",
    "2021-10-05 14:20:01.014883",
    "2021-10-05 14:20:01.014886",
    "2021-10-05 14:20:01.014888",
    "2021-10-05 14:20:01.014891",
    "2021-10-05 14:20:01.014893",
    "2021-10-05 14:20:01.014896",
    "2021-10-05 14:20:01.014899",
    "2021-10-05 14:20:01.014902",
    "2021-10-05 14:20:01.014904",
    "2021-10-05 14:20:01.014907",
    "2021-10-05 14:20:01.014909",
    "2021-10-05 14:20:01.014912",
    "2021-10-05 14:20:01.014914",
    "2021-10-05 14:20:01.014917",
    "2021-10-05 14:20:01.014920",
    "2021-10-05 14:20:01.014922",
    "2021-10-05 14:20:01.014925",
    "2021-10-05 14:20:01.014927",
    "2021-10-05 14:20:01.014931",
    "2021-10-05 14:20:01.014934",
    "2021-10-05 14:20:01.014936",
    "2021-10-05 14:20:01.014939",
    "2021-10-05 14:20:01.014942",
    "2021-10-05 14:20:01.014944",
    "2021-10-05 14:20:01.014947",
    "2021-10-05 14:20:01.014950",
    "2021-10-05 14:20:01.014952",
    "2021-10-05 14:20:01.014955",
    "2021-10-05 14:20:01.014957",
    "2021-10-05 14:20:01.014960",
    "2021-10-05 14:20:01.014964",
    "2021-10-05 14:20:01.014967",
    "2021-10-05 14:20:01.014970",
    "2021-10-05 14:20:01.014973",
    "2021-10-05 14:20:01.014975",
    "2021-10-05 14:20:01.014978",
    "2021-10-05 14:20:01.014980",
    "2021-10-05 14:20:01.014983",
    "2021-10-05 14:20:01.014986"
)
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add the metadata sheet after the data sheet
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Copy header style (bold/border/center) and index-column style from the data sheet
$dataSheet.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Dilated Cardiomyopathy and conduction defects"
$newSheet.Range("C2").Value = 47
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.70"
$newSheet.Range("E2").Value = "2021-05-12T14:07:44.532468Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:01.011491"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/47/?format=json"

$dataSheet.Select()
Write-Host "done"
